# Update elapsed-duration values in the Active_Outages.xlsx workbook.
# These "Elapsed Duration(Hrs)" cells (column G) are stored as plain text
# (e.g. "3883:34:46") and were recalculated to newer elapsed-time values
# (or cleared, in the one case where the underlying record became empty).

$wb = $excel.ActiveWorkbook

# Sheet "R1"
$ws = $wb.Worksheets.Item("R1")
$ws.Range("G2").Value = "3918:31:03"
$ws.Range("G3").Value = "58:03:41"
$ws.Range("G4").Value = ""

# Sheet "R2"
$ws = $wb.Worksheets.Item("R2")
$ws.Range("G2").Value = "12099:53:28"
$ws.Range("G3").Value = "3229:36:57"
$ws.Range("G4").Value = "467:48:31"

# Sheet "R4"
$ws = $wb.Worksheets.Item("R4")
$ws.Range("G2").Value = "2945:43:17"
$ws.Range("G3").Value = "172:55:32"

# Sheet "R5"
$ws = $wb.Worksheets.Item("R5")
$ws.Range("G2").Value = "419:42:16"

# Sheet "R6"
$ws = $wb.Worksheets.Item("R6")
$ws.Range("G2").Value = "60:14:34"
